$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "No "
$ws.Range("B1").Value = "Question "
$ws.Range("C1").Value = "Answer"

# Answer column values - set C3 before C2 so the shared-string table order
# matches the source workbook (longer "This file defines..." text first,
# then the stm32f7xx_hal_msp.c text)
$ws.Range("C3").Value = "This file defines all initialization functions to configure the IP instances according to the user configuration ( pin allocation, enabling of lock, use of DMA and interupt)"
$ws.Range("C2").Value = """stm32f7xx_hal_msp.c"" (MSP ->MCU Support Package ), t"

# New column width for the Answer column
$ws.Columns.Item(3).ColumnWidth = 18.498697916666668

# Wrap text for the question/answer body cells
$ws.Range("B2:C5").WrapText = $true

# Row heights for the wrapped long-answer rows
$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 115.2

# Selection, matching the committed state
$ws.Range("C8").Select() | Out-Null
